$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-8 from 45174 to 45175
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45175
}
